# Update the "Förändrad" (Changed) date column (C) for rows 2-21
# from 2023-10-09 (serial 45208) to 2023-10-13 (serial 45212).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2:C21").Value = 45212
